$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume snapshot data (GitHub Actions scheduled update).
# The "Price" column (D) sometimes holds plain decimal-looking text (e.g. "228.60")
# which Excel would otherwise silently coerce into a real number when assigned
# directly. Each such cell is temporarily switched to Text format ("@") before its
# value is written, then restored to the default "Normal" style afterwards so the
# cell formatting ends up exactly as it started (only the displayed text changes).
$priceCells = @("D2","D3","D5","D7","D10","D12","D13","D14","D17","D18","D19","D20","D21","D22","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D38","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '39.358.81'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '2.158.67'
$ws.Range("E3").Value = '  +3.10%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '228.60'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("D7").Value = '64.25'
$ws.Range("E7").Value = '  +5.23%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +2.78%  '
$ws.Range("D10").Value = '0.0857'
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '15.92'
$ws.Range("E12").Value = '  +3.31%  '
$ws.Range("D13").Value = '2.481.03'
$ws.Range("E13").Value = '  -19.11%  '
$ws.Range("D14").Value = '22.26'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '2.153.56'
$ws.Range("E17").Value = '  +3.21%  '
$ws.Range("D18").Value = '39.278.97'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '71.94'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '6.14'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '0.0₃0852'
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("D22").Value = '231.33'
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("E24").Value = '  +5.59%  '
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("D26").Value = '9.63'
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("D27").Value = '172.07'
$ws.Range("D28").Value = '0.139'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").Value = '19.93'
$ws.Range("E29").Value = '  +3.31%  '
$ws.Range("D30").Value = '1.41'
$ws.Range("E30").Value = '  -1.81%  '
$ws.Range("D31").Value = '2.70'
$ws.Range("E31").Value = '  +8.67%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").Value = '4.62'
$ws.Range("E33").Value = '  +2.48%  '
$ws.Range("D34").Value = '4.78'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").Value = '7.09'
$ws.Range("E35").Value = '  +7.60%  '
$ws.Range("D36").Value = '0.0618'
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("D38").Value = '3.59'
$ws.Range("E38").Value = '  -0.49%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = '103.94'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").Value = '0.0230'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '17.82'
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").Value = '1.539.78'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  +4.07%  '
$ws.Range("D45").Value = '4.30'
$ws.Range("E45").Value = '  +4.30%  '
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").Value = '2.82'
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.0924'
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").Value = '  +5.51%  '
$ws.Range("D49").Value = '7.78'
$ws.Range("E49").Value = '  +1.59%  '
$ws.Range("D50").Value = '2.363.82'
$ws.Range("E50").Value = '  +3.18%  '
$ws.Range("D51").Value = '2.97'
$ws.Range("E51").Value = '  -0.30%  '

foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}

